$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32 (shifting the existing Small_intestine / Large_intestine
# rows down by one) to make room for a new Brain (dose=30) data point.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row with the new Brain / dose=30 / %ID=0.01 datapoint.
# (The inserted row already inherits the surrounding A/B/D/E formatting automatically.)
$ws.Cells.Item(32, 1).Value = "Brain"
$ws.Cells.Item(32, 2).Value = 30
$ws.Cells.Item(32, 3).Value = 0.01
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = "iv"

# Column C gets a distinct style here: number format 0.00, bold font, centered.
$c32 = $ws.Cells.Item(32, 3)
$c32.NumberFormat = "0.00"
$c32.Font.Bold = $true
$c32.HorizontalAlignment = -4108

# Restore the view state recorded after the edit.
$ws.Application.Goto($ws.Range("A11"), $true)
$ws.Range("H29").Select()
